$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegisterData")

$ws.Range("A4").Value = "situ"
$ws.Range("B4").Value = "kumar"
$ws.Range("C4").Value = "Babamama143"
$ws.Range("D4").Value = "Babamama143"

$ws.Range("D4").Select()
